$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2,3) {
    $ws.Range("G$row").Value = 0.1180102915951973
    $ws.Range("H$row").Value = 0.1180102915951973
    $ws.Range("I$row").Value = 0.1080617495711835
    $ws.Range("J$row").Value = 0.09409436522362459
    $ws.Range("K$row").Value = 28.8
    $ws.Range("L$row").Value = 0.09879931389365351
    $ws.Range("M$row").Value = 6.9
    $ws.Range("N$row").Value = 0.02436440677966102
    $ws.Range("O$row").Value = 0.2395833333333333
    $ws.Range("P$row").Value = 6.9
    $ws.Range("Q$row").Value = 0.02436440677966102
    $ws.Range("R$row").Value = 0.2395833333333333
    $ws.Range("U$row").Value = 36
    $ws.Range("V$row").Value = 0.1271186440677966
    $ws.Range("W$row").Value = 0.2352941176470588
    $ws.Range("X$row").Value = 0.06375921118516384
    $ws.Range("Y$row").Value = 0.171534906461895
    $ws.Range("Z$row").Value = 2.89992041384799
    $ws.Range("AA$row").Value = 0.2728661705400574
    $ws.Range("AB$row").Value = 0.06370550425146632
    $ws.Range("AC$row").Value = 0.209160666288591
    $ws.Range("AD$row").Value = 1.36
    $ws.Range("AF$row").Value = 1.36
    $ws.Range("AG$row").Value = -34.64
    $ws.Range("AH$row").Value = 0.004779308405960079
    $ws.Range("AI$row").Value = 0.009433962264150943
    $ws.Range("AJ$row").Value = -0.1393627293208883
    $ws.Range("AK$row").Value = -0.3202662721893491
    $ws.Range("AL$row").Value = 0.111
    $ws.Range("AM$row").Value = 0.111
    $ws.Range("AN$row").Value = 0.03976608187134503
    $ws.Range("AO$row").Value = 283.7837837837838
    $ws.Range("AP$row").Value = -1.012865497076023
    $ws.Range("AQ$row").Value = 283.7837837837838
}
